# [File] FileService 구현 및 테스트 완료
#
# - "파일" sheet (FileService 관련domain sheet): mark the "Resource" /
#   "Base64" read-DTO rows (rows 4 & 5, columns E/F = 구현/테스트 여부) as
#   done ("O"), now that FileService implementation + tests are complete.
# - Update the active sheet/selection bookkeeping left behind by the edit:
#   the author ended up on the "파일" tab at F6, having last clicked B13 on
#   "도메인".

$wb = $excel.ActiveWorkbook

# --- 파일 (2nd sheet): fill in "구현"/"테스트" (E/F) = "O" for rows 4 & 5 ---
$wsFile = $wb.Worksheets.Item(2)
$wsFile.Range("E4").Value = "O"
$wsFile.Range("F4").Value = "O"
$wsFile.Range("E5").Value = "O"
$wsFile.Range("F5").Value = "O"

# --- 도메인 (1st sheet): was the selected tab before, selection moved to B13 ---
$wsDomain = $wb.Worksheets.Item(1)
$wsDomain.Activate() | Out-Null
$wsDomain.Range("B13").Select() | Out-Null

# --- 파일 becomes the active/selected tab, with selection on F6 ---
$wsFile.Activate() | Out-Null
$wsFile.Range("F6").Select() | Out-Null
